$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Excel auto-converts plain-numeric-looking text (e.g. "1.00", "582.62")
# into real numbers when assigned via .Value, which would silently drop
# formatting such as trailing zeros. Prefix those with an apostrophe so
# Excel stores them as literal text (exactly like typing '1.00 by hand).
function Set-TextValue($cell, [string]$text) {
    if ($text -match '^[+-]?[0-9]*\.?[0-9]+$') {
        $cell.Value = "'" + $text
    } else {
        $cell.Value = $text
    }
}

# row -> @(new D value or $null, new E value or $null)
$changes = @{
    2  = @("66.082.89", "  -0.28%  ")
    3  = @("3.327.55",  "  +0.56%  ")
    4  = @("0.998",     "  -0.25%  ")
    5  = @("582.62",    $null)
    6  = @("185.10",    "  -2.48%  ")
    7  = @("1.00",      "  -0.03%  ")
    8  = @("3.322.25",  "  +0.65%  ")
    9  = @("0.575",     "  -2.39%  ")
    10 = @($null,       "  -3.53%  ")
    11 = @($null,       "  -2.03%  ")
    12 = @("47.07",     "  -1.75%  ")
    13 = @($null,       "  -1.80%  ")
    14 = @("642.67",    "  +4.41%  ")
    15 = @("3.856.17",  "  +0.35%  ")
    16 = @($null,       "  -2.69%  ")
    17 = @("66.168.81", "  -0.32%  ")
    18 = @("17.95",     "  -0.88%  ")
    19 = @($null,       "  -0.06%  ")
    20 = @("3.326.21",  "  +0.28%  ")
    21 = @("11.06",     "  -0.70%  ")
    22 = @("0.897",     "  -1.66%  ")
    23 = @("17.90",     "  -3.34%  ")
    24 = @($null,       "  -1.31%  ")
    25 = @("100.36",    "  -1.47%  ")
    26 = @($null,       "  -0.73%  ")
    27 = @($null,       "  -0.23%  ")
    28 = @("9.50",      "  -2.62%  ")
    29 = @("31.20",     "  +2.80%  ")
    30 = @($null,       "  -2.18%  ")
    31 = @($null,       "  -1.07%  ")
    32 = @("595.18",    "  +2.91%  ")
    33 = @("3.85",      "  -6.04%  ")
    34 = @($null,       "  -1.43%  ")
    35 = @($null,       "  -0.43%  ")
    36 = @("3.840.52",  "  +2.53%  ")
    38 = @("55.91",     "  -2.38%  ")
    39 = @($null,       "  -4.75%  ")
    40 = @($null,       "  -3.06%  ")
    41 = @("2.66",      "  -3.13%  ")
    42 = @("32.73",     "  -4.17%  ")
    43 = @("3.43",      "  +4.62%  ")
    44 = @("3.17",      "  -5.45%  ")
    45 = @($null,       "  -1.90%  ")
    46 = @($null,       "  -3.81%  ")
    47 = @($null,       "  -13.37%  ")
    48 = @("0.127",     "  -1.54%  ")
    49 = @($null,       "  +0.17%  ")
    50 = @("2.55",      "  -2.08%  ")
    51 = @("130.61",    "  +6.68%  ")
}

foreach ($row in $changes.Keys) {
    $vals = $changes[$row]
    $dVal = $vals[0]
    $eVal = $vals[1]
    if ($null -ne $dVal) {
        Set-TextValue $ws.Cells.Item($row, 4) $dVal
    }
    if ($null -ne $eVal) {
        $ws.Cells.Item($row, 5).Value = $eVal
    }
}
